$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (F column) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9
$ws1.Range("F3").Value = 301
$ws1.Range("F5").Value = 2558
$ws1.Range("F6").Value = 1844
$ws1.Range("F7").Value = 354
$ws1.Range("F9").Value = 898

# Sheet "全部类型" (sheet4) - update "想去人数" (F column) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9
$ws4.Range("F3").Value = 301
$ws4.Range("F5").Value = 2558
$ws4.Range("F6").Value = 1844
$ws4.Range("F7").Value = 354
$ws4.Range("F10").Value = 898
